$d = $word.ActiveDocument

function Append-Done([string]$searchText) {
    $range = $d.Content
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Collapse(0)
        $range.InsertAfter(" - Done")
        $range.Font.Name = "Times New Roman"
        $range.Font.NameFarEast = "Times New Roman"
        $range.Font.NameBi = "Times New Roman"
        $range.Font.NameOther = "Times New Roman"
        $range.Font.Size = 12
    }
}

Append-Done("1.Modificarea(crearea) bazei de date")
Append-Done("2.Conexiunea cu baza de date (doctrine)")
